# Fix typo in slides (Assembly F'24)
#
# Slide with SlideID 738 has a shape (Id 4, "Rectangle 4") containing the
# assembly listing for mult2. The "subq" line reads:
#       subq  %rsp, 32           # allocate frame
# which is invalid AT&T syntax (operands reversed). It should read:
#       subq  $32, %rsp          # allocate frame
#
# We locate the slide/shape robustly (by SlideID / Shape.Id, not raw
# position), then rewrite just that one paragraph's runs in place so the
# surrounding run-level formatting (Courier New, sz=1800, b=0, etc.) is
# preserved.

$p = $ppt.ActivePresentation

$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq 738) {
        $targetSlide = $candidate
        break
    }
}

if ($targetSlide -eq $null) {
    throw "could not find slide with SlideID 738"
}

$targetShape = $null
for ($i = 1; $i -le $targetSlide.Shapes.Count; $i++) {
    $candidate = $targetSlide.Shapes.Item($i)
    if ($candidate.Id -eq 4) {
        $targetShape = $candidate
        break
    }
}

if ($targetShape -eq $null) {
    throw "could not find shape with Id 4 on slide 738"
}

$tr = $targetShape.TextFrame.TextRange

# Find the paragraph whose text is the "subq" allocate-frame line.
$targetParaIndex = -1
$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like '*subq*%rsp, 32*allocate frame*') {
        $targetParaIndex = $i
        break
    }
}

if ($targetParaIndex -eq -1) {
    throw "could not find the 'subq ... allocate frame' paragraph"
}

$para = $tr.Paragraphs($targetParaIndex, 1)

# Runs before the edit:
#   1 "  "
#   2 "subq"
#   3 "  %"
#   4 "rsp"
#   5 ", 32           "
#   6 "# allocate frame"
#
# Target runs after the edit:
#   1 "  "
#   2 "subq"
#   3 "  $32, "            (new run)
#   4 "%"
#   5 "rsp"
#   6 "          # allocate frame"

$run3 = $para.Runs(3, 1)
if ($run3.Text -ne '  %') {
    throw "unexpected run 3 text: [$($run3.Text)]"
}

# Step 1: split run 3 ("  %") right before its first character and plant
# the new "  $32, " text there. This turns run 3 into two runs: the new
# "  $32, " run (inheriting run 3's formatting) followed by the
# unmodified remainder " %".
$splitPoint = $para.Characters(7, 1)
$splitPoint.Text = '  $32, '

# Step 2: trim the remainder (now run 4, " %") down to just "%".
$para = $tr.Paragraphs($targetParaIndex, 1)
$remainder = $para.Runs(4, 1)
if ($remainder.Text -ne ' %') {
    throw "unexpected remainder run text: [$($remainder.Text)]"
}
$remainder.Text = '%'

# Step 3: delete the old ", 32           " run (now run 6).
$para = $tr.Paragraphs($targetParaIndex, 1)
$oldOperand = $para.Runs(6, 1)
if ($oldOperand.Text -ne ', 32           ') {
    throw "unexpected old-operand run text: [$($oldOperand.Text)]"
}
$oldOperand.Text = ''

# Step 4: pad "# allocate frame" (now run 6 again, after the deletion)
# with the leading spaces that used to belong to the deleted run, so the
# comment column stays aligned.
$para = $tr.Paragraphs($targetParaIndex, 1)
$comment = $para.Runs(6, 1)
if ($comment.Text -ne '# allocate frame') {
    throw "unexpected comment run text: [$($comment.Text)]"
}
$comment.Text = '          # allocate frame'
